$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 392, shifting existing rows 392:478 down to 393:479
$ws.Rows.Item(392).Insert()

# Populate the newly inserted row 392 with its values.
# Columns A,B,C,E,F,G,H,I,J,K are constant across the whole data set,
# copy them from the (now shifted) row 393 which still holds the old row-392 data.
$ws.Cells.Item(392, 1).Value = $ws.Cells.Item(393, 1).Value2
$ws.Cells.Item(392, 2).Value = $ws.Cells.Item(393, 2).Value2
$ws.Cells.Item(392, 3).Value = $ws.Cells.Item(393, 3).Value2
$ws.Cells.Item(392, 5).Value = $ws.Cells.Item(393, 5).Value2
$ws.Cells.Item(392, 6).Value = $ws.Cells.Item(393, 6).Value2
$ws.Cells.Item(392, 7).Value = $ws.Cells.Item(393, 7).Value2
$ws.Cells.Item(392, 8).Value = $ws.Cells.Item(393, 8).Value2
$ws.Cells.Item(392, 9).Value = $ws.Cells.Item(393, 9).Value2
$ws.Cells.Item(392, 10).Value = $ws.Cells.Item(393, 10).Value2
$ws.Cells.Item(392, 11).Value = $ws.Cells.Item(393, 11).Value2

$ws.Range("D392").Value = 44951
$ws.Range("L392").Value = "Primera"
$ws.Range("M392").Value = 220
$ws.Range("N392").Value = 7000
$ws.Range("O392").Value = 7500
$ws.Range("P392").Value = 7227
$ws.Range("Q392").Value = "`$/caja 7 kilos"
$ws.Range("R392").Value = "Región del Maule"
$ws.Range("S392").Value = 1032
$ws.Range("T392").Value = 7

# Match the date-formatted style used by the other D-column cells
$ws.Range("D392").NumberFormat = $ws.Range("D393").NumberFormat
